$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-5.56%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'2"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'40.70"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.19%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'2"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.76%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'2"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.07355"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.38%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'2"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'4.293"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.27%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'2"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'1.543"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-9.21%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'2"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.9237"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.44%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'2"
$ws.Range("G8").Style = "Normal"
$ws.Range("D9").Value = "'0.1200"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-3.11%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'2"
$ws.Range("G9").Style = "Normal"
$ws.Range("E10").Value = "'-2.88%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'2"
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D11").Value = "'0.04280"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.06%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'2"
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08598"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-5.21%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'2"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.1054"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.09%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'2"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.001276"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.40%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'2"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.005795"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.18%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'2"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'3.349"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.13%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'2"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'2.368"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.27%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'2"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'0.3286"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-2.07%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'2"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'7.676"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.12%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'2"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.1391"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'2.39%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'2"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.2884"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.95%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'2"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.03943"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'-1.75%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'2"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.001261"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.42%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'2"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'0.003779"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-7.25%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'2"
$ws.Range("G24").Style = "Normal"
$ws.Range("E25").Value = "'0.57%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'2"
$ws.Range("G25").Style = "Normal"
$ws.Range("D26").Value = "'0.0003726"
$ws.Range("D26").Style = "Normal"
$ws.Range("G26").Value = "'2"
$ws.Range("G26").Style = "Normal"
$ws.Range("G27").Value = "'2"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'2"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'2"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'2"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'2"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'2"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'2"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'2"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'2"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'2"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'2"
$ws.Range("G37").Style = "Normal"
$ws.Range("D38").Value = "'0.02318"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'-3.57%"
$ws.Range("E38").Style = "Normal"
$ws.Range("G38").Value = "'2"
$ws.Range("G38").Style = "Normal"
$ws.Range("D39").Value = "'0.04990"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-2.53%"
$ws.Range("E39").Style = "Normal"
$ws.Range("G39").Value = "'2"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.005646"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'119.61%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'2"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.007675"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-1.13%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'2"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.1284"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.44%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'2"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.007359"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.51%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'2"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.007090"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-11.59%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'2"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.3172"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.51%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'2"
$ws.Range("G45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006358"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-4.13%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'2"
$ws.Range("G46").Style = "Normal"
$ws.Range("E47").Value = "'-0.24%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'2"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.01862"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-93.13%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'2"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002102"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.24%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'2"
$ws.Range("G49").Style = "Normal"
$ws.Range("D50").Value = "'0.0002002"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.24%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'2"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'2"
$ws.Range("G51").Style = "Normal"
